$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Header row 8: period labels - shift left, add new rightmost period
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/10"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/10"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/10"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/10"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/10"

# Header row 9: publish dates - shift left, add new rightmost date
$ws.Range("D9").Value = "1399-02-31 (12)"
$ws.Range("E9").Value = "1400-02-27 (9)"
$ws.Range("F9").Value = "1401-02-25 (9)"
$ws.Range("G9").Value = "1402-02-24 (10)"
$ws.Range("H9").Value = "1402-02-30 (3)"

# Data rows 12-56: shift D<-E, E<-F, F<-G, G<-H, H<-new value
$data = @{
  12 = @("-", 191998, -6815, 362019, 1046786)
  13 = @(0, -9742, -15036, -44811, -66889)
  14 = @(178569, 182256, -21851, 317208, 979897)
  16 = @(5841, 0, 0, 0, 0)
  17 = @(-122400, -88809, -24867, -149064, -396870)
  18 = @("-", 0, 0, 0, 0)
  19 = @(0, 0, 0, 0, 0)
  20 = @(-214, 0, 0, 0, -12881)
  21 = @(0, 0, 0, 0, -67100)
  22 = @(0, 0, 0, 0, 50000)
  23 = @(0, 166925, 0, 0, 0)
  24 = @(0, 0, 0, 0, 0)
  25 = @(0, 0, 0, 0, 33)
  26 = @(0, 0, 0, 0, -26)
  27 = @(0, 0, 0, 0, 0)
  28 = @(0, 0, 0, 0, 0)
  29 = @("-", 0, 0, 0, 0)
  30 = @(0, 2781, 0, 121, 1393)
  31 = @(1048, 2098, 0, 0, 10363)
  32 = @(-116773, 82995, -24867, -148943, -415088)
  33 = @(5550, 265251, -46718, 168265, 564809)
  35 = @(0, 0, 145102, 0, 0)
  36 = @("-", "-", "-", "-", 0)
  37 = @(0, 0, 0, 0, 0)
  38 = @(0, 0, 0, 0, 0)
  39 = @(150987, 135000, 171000, 260000, 0)
  40 = @(-174141, -261059, -156281, -379444, -248094)
  41 = @(-57294, -121659, -26485, -98693, -167955)
  42 = @("-", 0, 0, 0, 0)
  43 = @("-", 0, 0, 0, 0)
  44 = @("-", 0, 0, 0, 0)
  45 = @("-", 0, 0, 0, 0)
  46 = @("-", 0, 0, 0, 0)
  47 = @("-", 0, 0, 0, 0)
  48 = @("-", 0, 0, 0, 0)
  49 = @("-", 0, 0, 0, 0)
  50 = @(0, 0, 0, 0, -259)
  51 = @(-23154, -247718, 133336, -218137, -416308)
  52 = @(-17604, 17533, 86618, -49872, 148501)
  53 = @(27183, 9579, 27112, 116086, 59557)
  54 = @(0, 0, 2356, -6657, 103307)
  55 = @(9579, 27112, 116086, 59557, 311365)
  56 = @(0, 0, 642398, 8638, 0)
}

$cols = @("D", "E", "F", "G", "H")
foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt 5; $i++) {
    $addr = "$($cols[$i])$r"
    $ws.Range($addr).Value = $vals[$i]
  }
}
